# Model training (cont.) slide: add a bullet point about the classification
# report being generated from the 20% held-out test set, punctuate the
# previous bullet, and grow the text box to fit the extra line (the shape
# uses spAutoFit, so PowerPoint grows it to match the new text extent).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(6)
$sh = $s.Shapes.Item(4)

$tr = $sh.TextFrame.TextRange

# Second bullet currently reads "...before optimization" (88 chars) starting
# right after the first bullet's 37 chars + 1 paragraph mark, i.e. at
# character 39. Add the missing trailing period without touching the first
# bullet's mixed (green) run formatting.
$para2 = $tr.Characters(39, 88)
$para2.Text = "Overall, a good performance as we have used small dataset as well as before optimization."

# Append a brand-new bullet (new paragraph) after the existing text, using
# the same bullet/line-spacing formatting inherited from the prior
# paragraph.
[void]$tr.InsertAfter("`rThis classification report is generated using 20% remaining test set. ")

# Grow the textbox height to match the autosized extent PowerPoint computes
# once the extra bullet is added (width stays the same).
$sh.Height = 232.4431
